# Refitting NCDEs to individual patients (for manuscript figure)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in column H (copy the header formatting from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Update the changed prediction/error values (rows 5, 6, 7 in the first block)
$ws.Range("D5").Value = 0.5940148994567231
$ws.Range("E5").Value = 0.5940148994567231

$ws.Range("D6").Value = 0.6321305949495737
$ws.Range("E6").Value = 0.6321305949495737

$ws.Range("F7").Value = 0.572897732257843

# New "Label" column values for block 1 (rows 2-7)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1

# New "Label" column values for block 2 (rows 8-13)
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 1
